$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9213765263557434
$ws.Range("B1").Value = 1.297832131385803
$ws.Range("C1").Value = 2.259460926055908
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.724870443344116
